# Updated 2D training schedules, no break screen
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data values for rows 2-6, columns A-I (J stays the shared text "train_dim2_1")
$data = @(
    @(1, 0, 8, 4, 6, 4, -2, 54, 5),
    @(2, 1, 7, 6, 6, 5, -1, 65, 5),
    @(3, 1, 9, 2, 4, 1, -5, 21, 5),
    @(4, 0, 6, 3, 3, 3, -3, 43, 5),
    @(5, 3, 9, 5, 5, 2, -4, 32, 5)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $values = $data[$i]
    for ($c = 0; $c -lt $values.Length; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $values[$c]
    }
    $ws.Cells.Item($row, 10).Value = "train_dim2_1"
}

# Update the selection to I1, matching the new view state
$ws.Range("I1").Select()
